# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# "542599ee-4232-41c3-9af0-cc538acce734" row (row 3) on both the
# zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-19 22:34:22"
$zhcn.Range("H3").Value = "2016-03-19 22:34:41"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-19 22:34:25"
$dede.Range("H3").Value = "2016-03-19 22:34:47"
